# Generate Report for Handoff
# Updates the localization-status report: flips rows from "Handed back: in
# sync with en-US" to "Ready for handoff", bumps the associated timestamps,
# switches Priority from "ht" to "mt", records a new handoff timestamp +
# version-mismatch Error Detail for the f6fba0a3 file on the zh-cn sheet,
# and narrows/widens the Status/Error Detail columns to fit the new text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = "Ready for handoff"
$ovw.Range("F2").Value = "Ready for handoff"
$ovw.Range("G2").Value = "2016-09-06 16:13:25"
$ovw.Range("E3").Value = "Ready for handoff"
$ovw.Range("F3").Value = "Ready for handoff"
$ovw.Range("G3").Value = "2016-09-06 16:13:25"

$ovw.Columns.Item(5).ColumnWidth = 16.3
$ovw.Columns.Item(6).ColumnWidth = 16.3

# ---- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("H2").Value = "2016-09-06 16:13:08"
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-09-06 16:13:08"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f21ceaae36469c01310b0e55070389052fd75aeb/e2e/f6fba0a3-a8c9-4e05-93fb-594be715e3ac.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de586ae35827fa423df1e1a17b51fe66d44ee5fc/e2e/f6fba0a3-a8c9-4e05-93fb-594be715e3ac.md."

$zhcn.Columns.Item(3).ColumnWidth = 16.3
$zhcn.Columns.Item(16).ColumnWidth = 39.1

# ---- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "mt"
$dede.Range("H2").Value = "2016-09-06 16:13:25"
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-09-06 16:13:25"

$dede.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(16).ColumnWidth = 39.1
